$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Full target data (A=STORE_ID, B=PRODUCT_ID, C=P_NUM) for rows 2..22
$data = @(
    @(1,1,7),
    @(1,2,5),
    @(1,3,2),
    @(2,1,6),
    @(2,2,4),
    @(2,3,3),
    @(3,1,7),
    @(3,2,3),
    @(3,3,3),
    @(4,1,5),
    @(4,2,4),
    @(4,3,6),
    @(5,1,7),
    @(5,2,7),
    @(5,3,7),
    @(6,1,6),
    @(6,2,2),
    @(6,3,4),
    @(7,1,7),
    @(7,2,4),
    @(7,3,4)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}

$ws.Range("D20").Select()
